$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): columns C, D, E get new labels
$ws.Cells.Item(1, 3).Value = "prediction"
$ws.Cells.Item(1, 4).Value = "rejection-f"
$ws.Cells.Item(1, 5).Value = "max"

# Update data rows 2-12: column C gets the taxon text (same as column D),
# column D stays the same taxon text, column E becomes numeric 1
for ($r = 2; $r -le 12; $r++) {
    $taxon = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $taxon
    $ws.Cells.Item($r, 5).Value = 1
}
